$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = "Utility (Percent)"

for ($r = 3; $r -le 23; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2
    $jVal = $ws.Cells.Item($r, 10).Value2
    $kVal = $ws.Cells.Item($r, 11).Value2

    $ws.Cells.Item($r, 9).Value  = ($iVal.ToString() + " msec")
    $ws.Cells.Item($r, 10).Value = ($jVal.ToString() + " msec")
    $ws.Cells.Item($r, 11).Value = ($kVal.ToString() + " msec")
}
